$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update F22 and F23 from "nein" to "ja"
$ws.Range("F22").Value = "ja"
$ws.Range("F23").Value = "ja"

# Add new cell C24 with value "nein"
$ws.Range("C24").Value = "nein"
$ws.Range("C24").HorizontalAlignment = -4108  # xlCenter

# Update selection to F10
$ws.Range("F10").Select()
